# Update "Generate Report for Handback" timestamps.
# These cells hold plain text timestamps (not Excel dates), so we force
# text assignment to avoid Excel re-interpreting them as date serials.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 09:04:01"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 09:03:55"
$wsZhCn.Range("K2").Value = "2016-08-17 09:04:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 09:04:01"
$wsDeDe.Range("K2").Value = "2016-08-17 09:04:40"
